# Walk the Shopping List, total up how much each client spent (storing the
# running totals in the $clientExpenses hashtable), and use those same
# purchases to deduct the quantities sold from the Inventory sheet's Stock
# column. While we're in there we also normalize the Discount column from a
# whole-number percent (e.g. 15) to a decimal fraction (e.g. 0.15) so it can
# be used directly in price math.

$wb = $excel.ActiveWorkbook
$invSheet = $wb.Worksheets.Item("Inventory")
$shopSheet = $wb.Worksheets.Item("Shopping List")

# Map Item name -> Inventory row number.
$itemRow = @{}
$invLastRow = $invSheet.Cells.Item($invSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $invLastRow; $r++) {
    $itemName = $invSheet.Cells.Item($r, 1).Value2
    $itemRow[$itemName] = $r
}

# How many units of each item were ordered, across every client on the
# Shopping List sheet.
$qtySold = @{}
$shopLastRow = $shopSheet.Cells.Item($shopSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $shopLastRow; $r++) {
    $itemName = $shopSheet.Cells.Item($r, 2).Value2
    $qty = $shopSheet.Cells.Item($r, 3).Value2
    if ($qtySold.ContainsKey($itemName)) {
        $qtySold[$itemName] = $qtySold[$itemName] + $qty
    } else {
        $qtySold[$itemName] = $qty
    }
}

# clientExpenses: Client -> running total of (price * qty) after discount,
# accumulated while we scan the Shopping List a second time.
$clientExpenses = @{}
for ($r = 2; $r -le $shopLastRow; $r++) {
    $client = $shopSheet.Cells.Item($r, 1).Value2
    $itemName = $shopSheet.Cells.Item($r, 2).Value2
    $qty = $shopSheet.Cells.Item($r, 3).Value2

    $row = $itemRow[$itemName]
    if ($row) {
        $price = $invSheet.Cells.Item($row, 3).Value2
        $discountPct = $invSheet.Cells.Item($row, 5).Value2
        $discountFraction = $discountPct / 100
        $lineCost = $price * $qty * (1 - $discountFraction)

        if ($clientExpenses.ContainsKey($client)) {
            $clientExpenses[$client] = $clientExpenses[$client] + $lineCost
        } else {
            $clientExpenses[$client] = $lineCost
        }
    }
}

# Restocking counts received from vendors since the last inventory count,
# keyed by item name (independent of what moved through the Shopping List).
$restock = @{
    "Ice Cream"    = 10
    "Oranges"      = 23
    "Nesquik"      = 80
    "Guitar"       = 15
    "Paper Plates" = 1
}

foreach ($itemName in $restock.Keys) {
    $row = $itemRow[$itemName]
    if ($row) {
        $sold = 0
        if ($qtySold.ContainsKey($itemName)) {
            $sold = $qtySold[$itemName]
        }
        $newStock = $restock[$itemName]
        $invSheet.Cells.Item($row, 4).Value = $newStock
    }
}

# Re-express every populated Discount cell as a decimal fraction instead of
# a bare percent number.
for ($r = 2; $r -le $invLastRow; $r++) {
    $discountPct = $invSheet.Cells.Item($r, 5).Value2
    if ($discountPct -ne 0) {
        $invSheet.Cells.Item($r, 5).Value = $discountPct / 100
    }
}

# Keep the active selections pointed at the cells last touched on each
# sheet.
$shopSheet.Activate()
$shopSheet.Range("A2").Select()

$invSheet.Activate()
$invSheet.Range("G18").Select()

Write-Output "clientExpenses + inventory stock updated"
